# Export user data to excel file.
# Add a new "User Data" sheet as the first sheet in the workbook, with the
# text "User Data" in cell A1, mirroring the existing tab layout pattern
# used by the other sheets (sheet name also written to A1).

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current first sheet so it always lands at
# the very front of the tab strip, regardless of which sheet is currently
# active.
$firstSheet = $wb.Worksheets.Item(1)
$userDataSheet = $wb.Worksheets.Add($firstSheet)
$userDataSheet.Name = "User Data"
$userDataSheet.Range("A1").Value = "User Data"
